$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J4:AS4").Value = 1.016025614293143
$ws.Range("J5:AS5").Value = -0.1315854458144212
